$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 574.75
$ws.Range("J58").Value = 864.1539
$ws.Range("L58").Value = 2592.4617
$ws.Range("N58").Value = -2892.4617

$ws.Range("H64").Value = 71431496
$ws.Range("I64").Value = 125002000
$ws.Range("J64").Value = 4166.5
$ws.Range("K64").Value = 125002000
$ws.Range("L64").Value = 4166.5
$ws.Range("M64").Value = -125001752
$ws.Range("N64").Value = -4662.5

$ws.Range("H67").Value = 71431496
$ws.Range("I67").Value = 125002000
$ws.Range("J67").Value = 4166.5
$ws.Range("K67").Value = 125002000
$ws.Range("L67").Value = 4166.5
$ws.Range("M67").Value = -125001142
$ws.Range("N67").Value = -5882.5

$ws.Range("H116").Value = 3382
$ws.Range("I116").Value = 3382
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3382
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 60
$ws.Range("N116").ClearContents()

$ws.Range("H129").Value = 809.21875
$ws.Range("I129").Value = 633.5833
$ws.Range("J129").Value = 914.6
$ws.Range("K129").Value = 1900.7499
$ws.Range("L129").Value = 2743.8
$ws.Range("M129").Value = 3099.2501
$ws.Range("N129").Value = -12743.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 44168
$ws.Range("I6").Value = 49001.6
$ws.Range("J6").Value = 20000
$ws.Range("K6").Value = 49001.6
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = -48828.6
$ws.Range("N6").Value = -20346

$ws.Range("H32").Value = 8885.468999999999
$ws.Range("I32").Value = 7845.6055
$ws.Range("J32").Value = 18114.25
$ws.Range("K32").Value = 7845.6055
$ws.Range("L32").Value = 18114.25
$ws.Range("M32").Value = -7558.6055
$ws.Range("N32").Value = -18688.25

$ws.Range("H61").Value = 314147.3
$ws.Range("I61").Value = 1679.3103
$ws.Range("K61").Value = 1679.3103
$ws.Range("M61").Value = -1467.3103

$ws.Range("H113").Value = 34694
$ws.Range("J113").Value = 34694
$ws.Range("L113").Value = 34694
$ws.Range("N113").Value = -43372

$ws.Range("H136").Value = 314147.3
$ws.Range("I136").Value = 1679.3103
$ws.Range("K136").Value = 5037.9309
$ws.Range("M136").Value = -2487.9309

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 17813.5
$ws.Range("J95").Value = 17813.5
$ws.Range("L95").Value = 17813.5
$ws.Range("N95").Value = -23305.5

$ws.Range("H102").Value = 8655.444
$ws.Range("I102").Value = 3862.375
$ws.Range("K102").Value = 3862.375
$ws.Range("M102").Value = -617.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 27736.842
$ws.Range("J64").Value = 27736.842
$ws.Range("L64").Value = 27736.842
$ws.Range("N64").Value = -28232.842

$ws.Range("H67").Value = 27736.842
$ws.Range("J67").Value = 27736.842
$ws.Range("L67").Value = 27736.842
$ws.Range("N67").Value = -29452.842

$ws.Range("H68").Value = 19980
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 19980
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 19980
$ws.Range("N68").Value = -21478
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 19980
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 19980
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 59940
$ws.Range("N71").Value = -67428
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2322.3809
$ws.Range("I70").Value = 1605.8334
$ws.Range("J70").Value = 3277.7778
$ws.Range("K70").Value = 4817.5002
$ws.Range("L70").Value = 9833.3334
$ws.Range("M70").Value = -4502.5002
$ws.Range("N70").Value = -10463.3334

$ws.Range("H73").Value = 2322.3809
$ws.Range("I73").Value = 1605.8334
$ws.Range("J73").Value = 3277.7778
$ws.Range("K73").Value = 4817.5002
$ws.Range("L73").Value = 9833.3334
$ws.Range("M73").Value = -3725.5002
$ws.Range("N73").Value = -12017.3334

$ws.Range("H117").Value = 395
$ws.Range("I117").Value = 347.5
$ws.Range("J117").Value = 490
$ws.Range("K117").Value = 1042.5
$ws.Range("L117").Value = 1470
$ws.Range("M117").Value = 2399.5
$ws.Range("N117").Value = -8354

$ws.Range("H129").Value = 1031.4375
$ws.Range("J129").Value = 1486.1428
$ws.Range("L129").Value = 4458.428400000001
$ws.Range("N129").Value = -14458.4284

$ws.Range("H130").Value = 1500

$ws.Range("H131").Value = 905.5161000000001
$ws.Range("J131").Value = 1038.44
$ws.Range("L131").Value = 3115.32
$ws.Range("N131").Value = -13195.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 40451.43
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 40451.43
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 40451.43
$ws.Range("N130").Value = -50491.43
$ws.Range("M130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value = 36000
$ws.Range("J97").Value = 36000
$ws.Range("L97").Value = 36000
$ws.Range("N97").Value = -37982

$ws.Range("H130").Value = 39750
$ws.Range("J130").Value = 39750
$ws.Range("L130").Value = 39750
$ws.Range("N130").Value = -49790

$ws.Range("H136").Value = 6353.923
$ws.Range("I136").Value = 1869.6923
$ws.Range("J136").Value = 10838.154
$ws.Range("K136").Value = 5609.0769
$ws.Range("L136").Value = 32514.462
$ws.Range("M136").Value = -3059.0769
$ws.Range("N136").Value = -37614.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 4325
$ws.Range("I8").Value = 1500
$ws.Range("J8").Value = 5266.6665
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 5266.6665
$ws.Range("M8").Value = -1360
$ws.Range("N8").Value = -5546.6665

$ws.Range("H57").Value = 35398.332
$ws.Range("J57").Value = 35398.332
$ws.Range("L57").Value = 35398.332
$ws.Range("N57").Value = -36906.332

$ws.Range("H76").Value = 30915.334
$ws.Range("J76").Value = 30915.334
$ws.Range("L76").Value = 30915.334
$ws.Range("N76").Value = -31545.334

$ws.Range("H79").Value = 30915.334
$ws.Range("J79").Value = 30915.334
$ws.Range("L79").Value = 30915.334
$ws.Range("N79").Value = -33099.334

$ws.Range("H107").Value = 342.84616
$ws.Range("I107").Value = 306.66666
$ws.Range("J107").Value = 777
$ws.Range("K107").Value = 919.9999799999999
$ws.Range("L107").Value = 2331
$ws.Range("M107").Value = 1000.00002
$ws.Range("N107").Value = -6171

$ws.Range("H136").Value = 6694.8237
$ws.Range("I136").Value = 10510.2
$ws.Range("J136").Value = 1244.2858
$ws.Range("K136").Value = 31530.6
$ws.Range("L136").Value = 3732.8574
$ws.Range("M136").Value = -28980.6
